# "bill model in separate file"
# Update the headers/data on the existing sheet and split the "Robert & fils"
# bill off into its own new worksheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Columns D (start date), E (end date) and F (price) hold plain text in this
# workbook (e.g. "04.11.2019", "300"), not real dates/numbers - force text
# formatting first so Excel doesn't auto-convert them.
$ws1.Range("D1:F4").NumberFormat = "@"

# --- Update headers (row 1) to French labels ---
$ws1.Range("A1").Value = "Type de traveaux"
$ws1.Range("B1").Value = "Nom de l'entreprise"
$ws1.Range("C1").Value = "Commentaire"
$ws1.Range("D1").Value = "Date de debut des traveaux"
$ws1.Range("E1").Value = "Date de fin des traveaux"
$ws1.Range("F1").Value = "Prix"
$ws1.Range("G1").Value = "Etat du payment"

# --- Update row 2 (keep A2/B2, rework the rest) ---
$ws1.Range("C2").Value = "Remplacer le robinet au sous sol du B3. Travail effectué, facture payé."
$ws1.Range("D2").Value = "04.11.2019"
$ws1.Range("E2").Value = "04.11.2019"
$ws1.Range("F2").Value = "300"
$ws1.Range("G2").Value = "Payé"

# --- Update row 3 (keep A3/B3, rework the rest) ---
$ws1.Range("C3").Value = "test"
$ws1.Range("D3").Value = "04.11.2019"
$ws1.Range("E3").Value = "04.11.2019"
$ws1.Range("F3").Value = "365"
$ws1.Range("G3").Value = "Payé"

# --- Add new row 4 ---
$ws1.Range("A4").Value = "plomberie"
$ws1.Range("B4").Value = "je suis sous l'eau"
$ws1.Range("C4").Value = "test1"
$ws1.Range("D4").Value = "04.11.2019"
$ws1.Range("E4").Value = "04.11.2019"
$ws1.Range("F4").Value = "60"
$ws1.Range("G4").Value = "Payé"

# --- Add a new worksheet "Robert & fils" after the existing sheet ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Robert & fils"

$ws2.Range("D1:F2").NumberFormat = "@"

$ws2.Range("A1").Value = "Type de traveaux"
$ws2.Range("B1").Value = "Nom de l'entreprise"
$ws2.Range("C1").Value = "Commentaire"
$ws2.Range("D1").Value = "Date de debut des traveaux"
$ws2.Range("E1").Value = "Date de fin des traveaux"
$ws2.Range("F1").Value = "Prix"
$ws2.Range("G1").Value = "Etat du payment"

$ws2.Range("A2").Value = "plomberie"
$ws2.Range("B2").Value = "Robert & fils"
$ws2.Range("C2").Value = "Changement du radiateur dans la loge de gardien. Fait, payé."
$ws2.Range("D2").Value = "04.11.2019"
$ws2.Range("E2").Value = "10.11.2019"
$ws2.Range("F2").Value = "1000"
$ws2.Range("G2").Value = "Payé"
